$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1423.1034
$ws.Range("I127").Value = 435
$ws.Range("J127").Value = 1799.5238
$ws.Range("K127").Value = 1305
$ws.Range("L127").Value = 5398.5714
$ws.Range("M127").Value = 3655
$ws.Range("N127").Value = -15318.5714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20841662
$ws.Range("I2").Value = 27788216
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 27788216
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -27788103
$ws.Range("N2").Value = -2226

$ws.Range("H5").Value = 62.5
$ws.Range("I5").Value = 72.666664
$ws.Range("J5").Value = 32
$ws.Range("K5").Value = 72.666664
$ws.Range("L5").Value = 32
$ws.Range("M5").Value = 39.333336
$ws.Range("N5").Value = -256

$ws.Range("H32").Value = 12248.051
$ws.Range("I32").Value = 8786.379000000001
$ws.Range("J32").Value = 39626.727
$ws.Range("K32").Value = 8786.379000000001
$ws.Range("L32").Value = 39626.727
$ws.Range("M32").Value = -8499.379000000001
$ws.Range("N32").Value = -40200.727

$ws.Range("H45").Value = 33340944
$ws.Range("I45").Value = 45464350
$ws.Range("K45").Value = 45464350
$ws.Range("M45").Value = -45463973

$ws.Range("H74").Value = 6251264.5
$ws.Range("I74").Value = 7576831
$ws.Range("J74").Value = 2164.2144
$ws.Range("K74").Value = 7576831
$ws.Range("L74").Value = 2164.2144
$ws.Range("M74").Value = -7575957
$ws.Range("N74").Value = -3912.2144

$ws.Range("H77").Value = 6251264.5
$ws.Range("I77").Value = 7576831
$ws.Range("J77").Value = 2164.2144
$ws.Range("K77").Value = 37884155
$ws.Range("L77").Value = 10821.072
$ws.Range("M77").Value = -37879787
$ws.Range("N77").Value = -19557.072

$ws.Range("H110").Value = 1479.091
$ws.Range("I110").Value = 653.9231
$ws.Range("J110").Value = 2671
$ws.Range("K110").Value = 653.9231
$ws.Range("L110").Value = 2671
$ws.Range("M110").Value = 1391.0769
$ws.Range("N110").Value = -6761

$ws.Range("H116").Value = 20841662
$ws.Range("I116").Value = 27788216
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 27788216
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = -27785922
$ws.Range("N116").Value = -6588

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20841662
$ws.Range("I3").Value = 27788216
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 27788216
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -27788102
$ws.Range("N3").Value = -2228

$ws.Range("H4").Value = 62.5
$ws.Range("I4").Value = 72.666664
$ws.Range("J4").Value = 32
$ws.Range("K4").Value = 72.666664
$ws.Range("L4").Value = 32
$ws.Range("M4").Value = 42.333336
$ws.Range("N4").Value = -262

$ws.Range("H22").Value = 343.125
$ws.Range("I22").Value = 277.7143
$ws.Range("J22").Value = 801
$ws.Range("K22").Value = 277.7143
$ws.Range("L22").Value = 801
$ws.Range("M22").Value = -104.7143
$ws.Range("N22").Value = -1147

$ws.Range("H94").Value = 1082.4147
$ws.Range("I94").Value = 730.34485
$ws.Range("J94").Value = 1933.25
$ws.Range("K94").Value = 730.34485
$ws.Range("L94").Value = 1933.25
$ws.Range("M94").Value = -279.34485
$ws.Range("N94").Value = -2835.25

$ws.Range("H105").Value = 1895732.6
$ws.Range("I105").Value = 2274379.2
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2274379.2
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -2272632.2
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 87.39130400000001
$ws.Range("I7").Value = 95.57143000000001
$ws.Range("J7").Value = 74.666664
$ws.Range("K7").Value = 95.57143000000001
$ws.Range("L7").Value = 74.666664
$ws.Range("M7").Value = 17.42856999999999
$ws.Range("N7").Value = -300.666664

$ws.Range("H22").Value = 859
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 916.55554
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 916.55554
$ws.Range("M22").Value = -250
$ws.Range("N22").Value = -1616.55554

$ws.Range("H31").Value = 11907675
$ws.Range("I31").Value = 20409222
$ws.Range("J31").Value = 5507.3145
$ws.Range("K31").Value = 20409222
$ws.Range("L31").Value = 5507.3145
$ws.Range("M31").Value = -20408927
$ws.Range("N31").Value = -6097.3145

$ws.Range("H34").Value = 11907675
$ws.Range("I34").Value = 20409222
$ws.Range("J34").Value = 5507.3145
$ws.Range("K34").Value = 20409222
$ws.Range("L34").Value = 5507.3145
$ws.Range("M34").Value = -20409020
$ws.Range("N34").Value = -5911.3145

$ws.Range("H62").Value = 2791.4707
$ws.Range("I62").Value = 2274.875
$ws.Range("J62").Value = 3250.6667
$ws.Range("K62").Value = 2274.875
$ws.Range("L62").Value = 3250.6667
$ws.Range("M62").Value = -1650.875
$ws.Range("N62").Value = -4498.6667

$ws.Range("H65").Value = 2791.4707
$ws.Range("I65").Value = 2274.875
$ws.Range("J65").Value = 3250.6667
$ws.Range("K65").Value = 11374.375
$ws.Range("L65").Value = 16253.3335
$ws.Range("M65").Value = -8254.375
$ws.Range("N65").Value = -22493.3335

$ws.Range("H132").Value = 1527.579
$ws.Range("I132").Value = 883.4074000000001
$ws.Range("J132").Value = 3108.7273
$ws.Range("K132").Value = 2650.2222
$ws.Range("L132").Value = 9326.1819
$ws.Range("M132").Value = -120.2222000000002
$ws.Range("N132").Value = -14386.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 420
$ws.Range("I7").Value = 84
$ws.Range("J7").Value = 980
$ws.Range("K7").Value = 252
$ws.Range("L7").Value = 2940
$ws.Range("M7").Value = -140
$ws.Range("N7").Value = -3164

$ws.Range("H129").Value = 1543.3334
$ws.Range("I129").Value = 1315
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 3945
$ws.Range("L129").Value = 6000
$ws.Range("M129").Value = 1055
$ws.Range("N129").Value = -16000

$ws.Range("H131").Value = 12245001
$ws.Range("I131").Value = 83333576
$ws.Range("J131").Value = 58388.027
$ws.Range("K131").Value = 250000728
$ws.Range("L131").Value = 175164.081
$ws.Range("M131").Value = -249995688
$ws.Range("N131").Value = -185244.081

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 85.28570999999999
$ws.Range("I2").Value = 91.333336
$ws.Range("J2").Value = 49
$ws.Range("K2").Value = 91.333336
$ws.Range("L2").Value = 49
$ws.Range("M2").Value = 21.666664
$ws.Range("N2").Value = -275

$ws.Range("H11").Value = 35457096
$ws.Range("I11").Value = 17500002
$ws.Range("J11").Value = 37933936
$ws.Range("K11").Value = 17500002
$ws.Range("L11").Value = 37933936
$ws.Range("M11").Value = -17499863
$ws.Range("N11").Value = -37934214

$ws.Range("H14").Value = 2563.125
$ws.Range("I14").Value = 3466.6667
$ws.Range("K14").Value = 3466.6667
$ws.Range("M14").Value = -3298.6667

$ws.Range("H24").Value = 447.72726
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 447.72726
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 447.72726
$ws.Range("N24").Value = -793.72726
$ws.Range("M24").ClearContents()

$ws.Range("H107").Value = 1360.826
$ws.Range("I107").Value = 1550.7142
$ws.Range("K107").Value = 1550.7142
$ws.Range("M107").Value = 369.2858000000001

$ws.Range("H132").Value = 5751222.5
$ws.Range("I132").Value = 7939816
$ws.Range("J132").Value = 6164.9375
$ws.Range("K132").Value = 23819448
$ws.Range("L132").Value = 18494.8125
$ws.Range("M132").Value = -23816918
$ws.Range("N132").Value = -23554.8125
